# Add 2022-Q1 fund-holdings sheet (between 2021-Q4 and 总计) and record it
# in the 总计 (totals) roll-up sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert the new "2022-Q1" worksheet right before the "总计" sheet.
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item(4)
$newSheet = $wb.Worksheets.Add($totalSheet)
$newSheet.Name = "2022-Q1"

# Match the page-margin convention used by the other data sheets (0.75in
# sides, 1in top/bottom, 0.5in header/footer == 54/54/72/72/36/36 pt).
$newSheet.PageSetup.LeftMargin = 54
$newSheet.PageSetup.RightMargin = 54
$newSheet.PageSetup.TopMargin = 72
$newSheet.PageSetup.BottomMargin = 72
$newSheet.PageSetup.HeaderMargin = 36
$newSheet.PageSetup.FooterMargin = 36

# ---------------------------------------------------------------------
# 2. Populate the header row (bold, centered, boxed - matches the other
#    quarterly sheets) and the single fund-holding data row.
# ---------------------------------------------------------------------
$headers = @{
    "B1" = "基金代码"
    "C1" = "基金名称"
    "D1" = "基金规模"
    "E1" = "股票总仓位"
    "F1" = "仓位占比"
    "G1" = "持有市值(亿元)"
    "H1" = "仓位排名"
}
foreach ($addr in $headers.Keys) {
    $cell = $newSheet.Range($addr)
    $cell.Value = $headers[$addr]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.Item(7).LineStyle = 1
    $cell.Borders.Item(8).LineStyle = 1
    $cell.Borders.Item(9).LineStyle = 1
    $cell.Borders.Item(10).LineStyle = 1
}

$rowIndex = $newSheet.Range("A2")
$rowIndex.Value = 0
$rowIndex.Font.Bold = $true
$rowIndex.HorizontalAlignment = -4108
$rowIndex.VerticalAlignment = -4160
$rowIndex.Borders.Item(7).LineStyle = 1
$rowIndex.Borders.Item(8).LineStyle = 1
$rowIndex.Borders.Item(9).LineStyle = 1
$rowIndex.Borders.Item(10).LineStyle = 1

function Set-TextValue($range, $value) {
    # Force the cell to stay text (so numeric-looking strings like "501069"
    # or "0.16" are not silently reinterpreted as numbers), then drop the
    # leftover text-number-format so no stray style sticks around.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $newSheet.Range("B2") "501069"
Set-TextValue $newSheet.Range("C2") "华宝标普中国Ａ股质量价值指数（ＬＯＦ）"
Set-TextValue $newSheet.Range("D2") "0.16"
Set-TextValue $newSheet.Range("E2") "94.73"
Set-TextValue $newSheet.Range("F2") "2.15"
Set-TextValue $newSheet.Range("G2") "0.0034"
$newSheet.Range("H2").Value = 9

# ---------------------------------------------------------------------
# 3. Update the "总计" sheet: insert a new top data row for 2022-Q1 and
#    renumber the existing rows' index column.
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows(2).Insert()

$a2 = $totalSheet.Range("A2")
$a2.Value = 0
$a2.Font.Bold = $true
$a2.HorizontalAlignment = -4108
$a2.VerticalAlignment = -4160
$a2.Borders.Item(7).LineStyle = 1
$a2.Borders.Item(8).LineStyle = 1
$a2.Borders.Item(9).LineStyle = 1
$a2.Borders.Item(10).LineStyle = 1

$totalSheet.Range("B2:D2").ClearFormats()
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 1
$totalSheet.Range("D2").Value = 0

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3

Write-Host "2022-Q1 sheet added and 总计 sheet updated"
